$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new rows of feed log data to the bottom of the sheet (rows 136 and 137)
$ws.Cells.Item(136, 1).Value = 135
$ws.Cells.Item(136, 2).Value = 1
$ws.Cells.Item(136, 3).Value = "2024-06-17 19:10:21"
$ws.Cells.Item(136, 4).Value = 200
$ws.Cells.Item(136, 5).Value = 14

$ws.Cells.Item(137, 1).Value = 136
$ws.Cells.Item(137, 2).Value = 2
$ws.Cells.Item(137, 3).Value = "2024-06-17 19:10:21"
$ws.Cells.Item(137, 4).Value = 200
$ws.Cells.Item(137, 5).Value = 1
